# Alex's updated dataset for 7.30.21
# Adds a "Number of Universities Worldwide in 2020 by Country" column (D)
# to the World Bank GDP dataset on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---
$ws.Range("D1").Value = "Number of Universities Worldwide in 2020 by Country"

# --- University counts for the countries that have data ---
$ws.Range("D2").Value  = 3254   # United States
$ws.Range("D3").Value  = 2595   # China
$ws.Range("D4").Value  = 1014   # Japan
$ws.Range("D5").Value  = 464    # Germany
$ws.Range("D6").Value  = 282    # United Kingdom
$ws.Range("D7").Value  = 4381   # India
$ws.Range("D8").Value  = 631    # France
$ws.Range("D9").Value  = 240    # Italy
$ws.Range("D10").Value = 370    # Canada
$ws.Range("D11").Value = 377    # Korea, Rep.
$ws.Range("D12").Value = 1096   # Russian Federation
$ws.Range("D13").Value = 1349   # Brazil
# Row 14 (Australia) intentionally has no university-count figure.
$ws.Range("D15").Value = 263    # Spain
$ws.Range("D16").Value = 1253   # Mexico
$ws.Range("D17").Value = 2694   # Indonesia

# --- Restore view state (active cell / scroll position) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 183
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C230").Select()
